# Update 0524 meeting discussion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "alpha 0.9"

# Header row: drop the "_0.1" suffix from the retrieval/validation column headers
$ws.Range("C1").Value = "檢索結果"
$ws.Range("D1").Value = "GPT_結果"
$ws.Range("E1").Value = "檢索驗證"
$ws.Range("F1").Value = "答案驗證"

# New column widths for A, B, E, F (C and D already have the right width)
$ws.Columns("A").ColumnWidth = 16.714285714285715
$ws.Columns("B").ColumnWidth = 17.714285714285715
$ws.Columns("E").ColumnWidth = 20.428571428571427
$ws.Columns("F").ColumnWidth = 15.0

# Move the active selection to C7
[void]$ws.Range("C7").Select()
